# "Generate Report for handoff"
# b.md is now ready for handoff: update the Overview + per-locale status
# sheets so b.md's status reads "Ready for handoff" and points at the
# freshly generated handoff files.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet ---------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-19 06:48:43"
foreach ($hl in $zh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# ---- de-de sheet ---------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$de.Range("D3").Value = "2016-01-19 06:48:53"
foreach ($hl in $de.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
